# Add a "DownloadDelta" column (timedelta of subprocess start/end) between
# "SizeLog" and "Last" in the SQLite export sheet.
#
# Before: Organization | Url | Folder | SizeWarc | SizeLog | Last       | State
# After : Organization | Url | Folder | SizeWarc | SizeLog | DownloadDelta | Last | State

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F is currently "Last". Insert a new blank column there, which
# shifts "Last" -> G and "State" -> H.
$ws.Columns.Item(6).Insert()

# Give the newly inserted column F its header.
$ws.Range("F1").Value = "DownloadDelta"

# Match the recorded selection left behind by the edit.
$ws.Range("F1").Select()
